# Update the "PRESUPUESTO" (column G) values on the "VENTA MENSUAL" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$updates = @{
    3  = 1500
    4  = 500
    5  = 1500
    6  = 250
    9  = 250
    12 = 1500
    13 = 1000
    14 = 2500
    16 = 1000
    17 = 3000
    18 = 3000
    19 = 1000
    20 = 1000
    21 = 2000
    26 = 500
    29 = 20500
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
